$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview" (1st sheet)
#   - row 2: Status columns B/C -> "Ready for handoff", D -> updated timestamp
#   - row 3 (f8d2db94-... entry): removed entirely
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item(1)

# The runtime only supports deleting *all* hyperlinks on a sheet in one shot,
# so remove them all up front and re-add just the ones that must survive.
$wsOverview.Hyperlinks.Delete()

$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-34-18 12:34:51"

$wsOverview.Rows.Item(3).Delete()

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/1206d0dd625b4af1eaccc09291ef603acb7abb54/e2e/73bb9739-d02c-4df7-9f3e-d5f82fd5f619.md", "", "", "73bb9739-d02c-4df7-9f3e-d5f82fd5f619.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (2nd sheet)
#   - row 2: Status (C) -> "Ready for handoff", Latest Handoff Datetime (E) updated
#   - row 3 (f8d2db94-... entry): removed entirely
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item(2)

$wsZhCn.Hyperlinks.Delete()

$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "2016-03-18 12:34:49"

$wsZhCn.Rows.Item(3).Delete()

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/1206d0dd625b4af1eaccc09291ef603acb7abb54/e2e/73bb9739-d02c-4df7-9f3e-d5f82fd5f619.md", "", "", "73bb9739-d02c-4df7-9f3e-d5f82fd5f619.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/1206d0dd625b4af1eaccc09291ef603acb7abb54/e2e/73bb9739-d02c-4df7-9f3e-d5f82fd5f619.md", "", "", ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1065caffdd41e7b42cd5deef5385136c863448b2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/73bb9739-d02c-4df7-9f3e-d5f82fd5f619.11ed010e27266a27c332c8fdd55168dead1be9f6.zh-cn.xlf", "", "", "73bb9739-d02c-4df7-9f3e-d5f82fd5f619.11ed010e27266a27c332c8fdd55168dead1be9f6.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/8c74938b634808059b7031b827d7ff4142c93548/e2e/73bb9739-d02c-4df7-9f3e-d5f82fd5f619.md", "", "", "73bb9739-d02c-4df7-9f3e-d5f82fd5f619.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/d889307e121d0bf1066bc21e8380eeffe09c9e94/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/73bb9739-d02c-4df7-9f3e-d5f82fd5f619.11ed010e27266a27c332c8fdd55168dead1be9f6.zh-cn.xlf", "", "", "73bb9739-d02c-4df7-9f3e-d5f82fd5f619.11ed010e27266a27c332c8fdd55168dead1be9f6.zh-cn.xlf")

# ---------------------------------------------------------------------------
# Sheet "de-de" (3rd sheet)
#   - row 2: Status (C) -> "Ready for handoff", Latest Handoff Datetime (E) updated
#   - row 3 (f8d2db94-... entry): removed entirely
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item(3)

$wsDeDe.Hyperlinks.Delete()

$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "2016-03-18 12:34:51"

$wsDeDe.Rows.Item(3).Delete()

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/1206d0dd625b4af1eaccc09291ef603acb7abb54/e2e/73bb9739-d02c-4df7-9f3e-d5f82fd5f619.md", "", "", "73bb9739-d02c-4df7-9f3e-d5f82fd5f619.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/1206d0dd625b4af1eaccc09291ef603acb7abb54/e2e/73bb9739-d02c-4df7-9f3e-d5f82fd5f619.md", "", "", ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/79a239f43b12ad0d01da36658e9c117831d4bac1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/73bb9739-d02c-4df7-9f3e-d5f82fd5f619.11ed010e27266a27c332c8fdd55168dead1be9f6.de-de.xlf", "", "", "73bb9739-d02c-4df7-9f3e-d5f82fd5f619.11ed010e27266a27c332c8fdd55168dead1be9f6.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/55c2fad09e87a0a175794fbfd6ac53c350949174/e2e/73bb9739-d02c-4df7-9f3e-d5f82fd5f619.md", "", "", "73bb9739-d02c-4df7-9f3e-d5f82fd5f619.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2d1adc8b07d58e127a4e90754f6ec42a50e2bc34/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/73bb9739-d02c-4df7-9f3e-d5f82fd5f619.11ed010e27266a27c332c8fdd55168dead1be9f6.de-de.xlf", "", "", "73bb9739-d02c-4df7-9f3e-d5f82fd5f619.11ed010e27266a27c332c8fdd55168dead1be9f6.de-de.xlf")

$wsOverview.Select()
